$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.555.57'
$ws.Range('E2').Value = '  -0.98%  '
$ws.Range('D3').Value = '1.884.77'
$ws.Range('E3').Value = '  +0.32%  '
$ws.Range('E4').Value = '  +0.29%  '
$ws.Range('D5').Value = "'326.46"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.36%  '
$ws.Range('E6').Value = '  +0.29%  '
$ws.Range('D7').Value = "'0.4585"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.76%  '
$ws.Range('D8').Value = "'0.3859"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.37%  '
$ws.Range('B9').Value = 'OKB'
$ws.Range('C9').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D9').Value = "'46.68"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.20%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').Value = "'0.07848"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.01%  '
$ws.Range('B11').Value = 'Polygon'
$ws.Range('C11').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D11').Value = "'0.9990"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.18%  '
$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D12').Value = "'21.60"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.77%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.901.69'
$ws.Range('E13').Value = '  +2.35%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').Value = "'7.053"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.30%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').Value = "'5.699"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.11%  '
$ws.Range('B16').Value = 'TRON'
$ws.Range('C16').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D16').Value = "'0.06963"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.38%  '
$ws.Range('B17').Value = 'Litecoin'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D17').Value = "'87.37"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.55%  '
$ws.Range('B18').Value = 'BinanceUSD'
$ws.Range('C18').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D18').Value = "'1.009"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.33%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = "'0.00001002"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.82%  '
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').Value = "'17.12"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.43%  '
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').Value = "'1.007"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.37%  '
$ws.Range('B22').Value = 'WrappedBTC'
$ws.Range('C22').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D22').Value = '28.583.65'
$ws.Range('E22').Value = '  -0.92%  '
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').Value = "'5.319"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.87%  '
$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D24').Value = "'10.96"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.58%  '
$ws.Range('B25').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C25').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D25').Value = '2.135.28'
$ws.Range('E25').Value = '  +3.02%  '
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').Value = "'2.061"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.78%  '
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').Value = "'154.60"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.55%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = "'19.35"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.33%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').Value = "'5.833"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.99%  '
$ws.Range('B30').Value = 'LidoDAOToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D30').Value = "'1.949"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.05%  '
$ws.Range('B31').Value = 'BitcoinCash'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D31').Value = "'118.15"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.45%  '
$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').Value = "'0.09319"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.94%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').Value = "'0.9208"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.42%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = "'5.295"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.71%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = "'1.332"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.69%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').Value = "'3.266"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.63%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').Value = "'0.05755"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.01%  '
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').Value = "'1.150"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.12%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = "'0.02063"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.28%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = "'7.843"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.27%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').Value = "'0.5666"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.30%  '
$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D42').Value = "'0.1787"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.55%  '
$ws.Range('B43').Value = 'Aptos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D43').Value = "'9.691"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.98%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = "'11.74"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.77%  '
$ws.Range('B45').Value = 'Cronos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D45').Value = "'0.07147"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.42%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = "'0.5342"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.21%  '
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').Value = "'2.174"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.48%  '
$ws.Range('D48').Value = "'1.118"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.61%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').Value = "'1.836"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.07%  '
$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D50').Value = "'112.36"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.79%  '
$ws.Range('B51').Value = 'MXToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D51').Value = "'2.469"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.22%  '
